# wzrv/expdata/1020.xlsx -- "Add files via upload"
#
# The sheet holds PHENIX W+/W- -> e asymmetry data points. This edit:
#  - retunes the row-2 data point (eta, stat_u)
#  - turns row 3 from a second "W+,Z" point into the "W-,Z" point, with its
#    own eta / stat_u / syst_u values (and the dependent A_L-error formula
#    recalculates from the new eta)
#  - drops the old row 4 (it was a duplicate "W+,Z" point) since the table
#    now only needs one row per boson charge
#  - bolds/centers the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: eta -> 0, stat_u -> 0.04 -------------------------------------
$ws.Range("D2").Value = 0
$ws.Range("I2").Value = 0.04

# --- Row 3: becomes the "W-,Z" point -------------------------------------
$ws.Range("G3").Value = "W-,Z"
$ws.Range("H3").Value = 0.17
$ws.Range("I3").Value = 0.08
$ws.Range("J3").Value = 0.02

# --- Drop the old row 4 (duplicate W+,Z point) ---------------------------
$ws.Rows("4:4").Delete()

# --- Header row: bold + centered -----------------------------------------
[void]($ws.Range("A1:M1").Font.Bold = $true)

# --- Restore the cursor to where the author left it ----------------------
[void]($ws.Range("K9").Select())
